$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Version: 5.0.0 -> 6.0.0
$ws.Cells.Item(3, 2).Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$ws.Cells.Item(8, 2).Value = "2022-01-21T20:46:54+00:00"

# Publisher value (was blank) -> "Alvearie Team"
$ws.Cells.Item(9, 2).Value = "Alvearie Team"

# Row 10 was "Contact" / "No display for ContactDetail" (duplicated on row 11).
# Replace row 10 with "Jurisdiction" / "United States of America" and delete
# the old duplicate row 11, shifting all following rows up by one.
$ws.Cells.Item(10, 1).Value = "Jurisdiction"
$ws.Cells.Item(10, 2).Value = "United States of America"
$ws.Rows.Item(11).Delete()

# Case Sensitive value (was blank) -> the text "true" (not a TRUE boolean).
# Typing "true" directly gets auto-coerced to a boolean by this engine, so
# build it as a formula result and paste back as a value to keep it a plain
# shared string - matches how the source workbook stores it.
$ws.Cells.Item(14, 2).Formula = '="true"'
$ws.Cells.Item(14, 2).Copy()
$ws.Cells.Item(14, 2).PasteSpecial(-4163)
